# Fix the typo "Comuter" -> "Computer" in the "College Name & Department"
# line of the title slide (slide 1), merging the three runs that made up
# "SNS COLLEGE OF TECHNOLOGY & " + "Comuter" + " Science And Engineering"
# into a single run reading
# "SNS COLLEGE OF TECHNOLOGY & Computer Science And Engineering".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 3")
$tr = $shp.TextFrame.TextRange

# Paragraph 3 is "College Name & Department : SNS COLLEGE OF TECHNOLOGY & Comuter Science And Engineering"
$para = $tr.Paragraphs(3)

# Underlying runs (by XML <a:r> order) in this paragraph are:
#   1: "College Name & Department : "
#   2: "SNS COLLEGE OF TECHNOLOGY & "
#   3: "Comuter"                         (flagged err="1")
#   4: " Science And Engineering"
# Clear out the trailing two runs first (from the end backwards) so the
# paragraph's overall text never balloons while we still hold references
# to them, then rewrite run 2 with the corrected, merged text.
$para.Runs(4).Text = ""
$para.Runs(3).Text = ""
$para.Runs(2).Text = "SNS COLLEGE OF TECHNOLOGY & Computer Science And Engineering"
